# Generate Report for Handoff
#
# This script updates the localization-status workbook to reflect that the
# "fd4b7196-3e80-496d-92b4-1397551f8b05" and "fd930e28-1aab-4b0e-96f4-a8b4bbb5b41d"
# records are now "Ready for handoff" (instead of "Handed back: in sync with
# en-US"), refreshes their handoff timestamps, records an error-detail message
# about a stale handback file, and widens the "Error Detail" column so the
# new message is readable.

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$newOverviewDate = "2016-08-27 16:26:44"

# The "Latest Handoff Datetime" refresh differs per target-language sheet.
$newHandoffDateBySheet = @{ "zh-cn" = "2016-08-27 16:26:39"; "de-de" = "2016-08-27 16:26:44" }

function Get-ErrorDetail($fileName) {
    return "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6892722a259a5b077180daaeb2404e3457160356/e2e/$fileName.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e81ffe3ef3aaecdfb4dad0e2115b840e0c876b0/e2e/$fileName.md."
}

# ---------------------------------------------------------------------------
# "Overview" sheet: rows 4 (fd4b7196...) and 5 (fd930e28...)
#   E/F (zh-cn / de-de status)            -> "Ready for handoff"
#   G   (Latest HO Xliff Generate Date)   -> "2016-08-27 16:26:44"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($row in 4..5) {
    $wsOverview.Cells.Item($row, 5).Value = $newStatus
    $wsOverview.Cells.Item($row, 6).Value = $newStatus
    $wsOverview.Cells.Item($row, 7).Value = $newOverviewDate
}

# ---------------------------------------------------------------------------
# "zh-cn" and "de-de" sheets: rows 4 (fd4b7196...) and 5 (fd930e28...)
#   C (Status)                   -> "Ready for handoff"
#   H (Latest Handoff Datetime)  -> refreshed timestamp
#   P (Error Detail)             -> stale-handback-file message
#   column P width                -> 40
# ---------------------------------------------------------------------------
$fileNames = @{ 4 = "fd4b7196-3e80-496d-92b4-1397551f8b05"; 5 = "fd930e28-1aab-4b0e-96f4-a8b4bbb5b41d" }

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $handoffDate = $newHandoffDateBySheet[$sheetName]

    foreach ($row in 4..5) {
        $ws.Cells.Item($row, 3).Value = $newStatus
        $ws.Cells.Item($row, 8).Value = $handoffDate
        $ws.Cells.Item($row, 16).Value = Get-ErrorDetail($fileNames[$row])
    }

    $ws.Columns("P").ColumnWidth = 39.17
}
